# Reorders the comma-separated "Recorded By" audit entries in column G so that
# the first-listed actor is moved to the end of the list (left-rotation by one),
# e.g. "dnasr281@gmail.com, System" -> "System, dnasr281@gmail.com".
# The one exact value "admin@admin.com, System" is left untouched, matching the
# source diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$usedRange = $ws.UsedRange
$firstRow = $usedRange.Row
$lastRow = $firstRow + $usedRange.Rows.Count - 1

# Locate the "Recorded By" column by header text on the first row, falling back
# to column G (7) if it cannot be found.
$col = 7
for ($c = 1; $c -le $usedRange.Columns.Count; $c++) {
    $header = $ws.Cells.Item($firstRow, $c).Value2
    if ($header -eq "Recorded By") {
        $col = $c
        break
    }
}

for ($r = $firstRow + 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $col)
    $val = $cell.Value2

    if ($null -eq $val) {
        continue
    }

    if ($val -eq "admin@admin.com, System") {
        continue
    }

    $parts = $val -split ", "
    if ($parts.Length -lt 2) {
        continue
    }

    $rotated = ($parts[1..($parts.Length - 1)] + $parts[0]) -join ", "
    $cell.Value = $rotated
}
